# "big modifications for elliminating Inputs"
# Duplicate the "P1" sheet twice, inserting the copies ("P2" and "P3")
# between "P1" and "Resources", then leave "P2" as the active/selected tab
# (with its own selection), matching the target diff.

$wb = $excel.ActiveWorkbook

$p1 = $wb.Worksheets.Item("P1")

# First duplicate of P1 -> placed right after P1, renamed "P2"
[void]$p1.Copy($null, $p1)
$wb.Worksheets.Item(2).Name = "P2"
$p2 = $wb.Worksheets.Item("P2")

# Second duplicate of P1 -> placed right after P2, renamed "P3"
[void]$p1.Copy($null, $p2)
$wb.Worksheets.Item(3).Name = "P3"

# Make "P2" the active sheet/tab, with its own selection (B21)
$p2.Activate()
[void]$p2.Range("B21").Select()
